# recebimentos e sangrias implementados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date update
$ws.Range("A1").Value = "CAIXA DIA 15/11/2025"

# Row 4 - pacote doce
$ws.Range("B4").Value = 33
$ws.Range("C4").Value = 196.38

# Row 5 - pacote sal
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 100.37

# Row 6 - pacote suiça
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 101.45

# Row 7 - pacote nazare
$ws.Range("B7").Value = 26
$ws.Range("C7").Value = 142.63
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 264

# Row 8 - pacote queijo
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 57.84

# Row 9 - pacote coquinho
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 535.23

# Row 10 - sequilho
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 14.96
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 388.98

# Row 11 - fardo
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 39.7

# Row 12 - recheado
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 40

# Row 13 - rosquinha
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 40

# Row 15 - bolo de rolo grande
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 10

# Row 16 - bolo de rolo pequeno
$ws.Range("B16").Value = $null
$ws.Range("C16").Value = $null

# Row 17 - doce de leite pequeno
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = $null

# Row 21 - nego bom
$ws.Range("B21").Value = $null
$ws.Range("C21").Value = $null

# Row 23 - beiju recheado
$ws.Range("B23").Value = $null
$ws.Range("C23").Value = $null

# Row 24 - biscoito alexandre
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 15

# Row 25 - brigadeiro
$ws.Range("B25").Value = 6
$ws.Range("C25").Value = 21

# Row 26 - pingo bel
$ws.Range("B26").Value = $null
$ws.Range("C26").Value = $null

# Row 30 - cocada grande
$ws.Range("B30").Value = $null
$ws.Range("C30").Value = $null

# Row 31 - cocada pequena
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 3.5

# Row 34 - coquero
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 3

# Row 36 - coca
$ws.Range("B36").Value = $null
$ws.Range("C36").Value = $null

# Row 42 - agua
$ws.Range("B42").Value = 3
$ws.Range("C42").Value = 8.960000000000001

# Row 43 - agua c/ gas
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 4

# Row 51 - OVOS (despesas)
$ws.Range("A51").Value = "OVOS"
$ws.Range("B51").Value = "24,00"

# Row 52 - VASSOURA (despesas)
$ws.Range("A52").Value = "VASSOURA"
$ws.Range("B52").Value = "30,00"

# Row 60 - DINHEIRO
$ws.Range("B60").Value = "183,00"

# Row 61 - DEBITO
$ws.Range("B61").Value = "182,00"

# Row 62 - PIX
$ws.Range("B62").Value = "328,00"
